# SR-9373:R4, R7 — adding test data edits to the Cora Intake Manifest sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cora Intake Manifest")

# Specimen Type / Specimen Source changed for every data row (7-11)
$ws.Range("G7:G11").Value = "Fresh Bone Marrow"
$ws.Range("H7:H11").Value = "Bone Marrow"

# Row 11 (5th sample) picked up a handful of additional corrections
$ws.Range("K11").Value = "3/24/2022"
$ws.Range("AC11").Value = 101.1
$ws.Range("AE11").Value = "12/16/2020"
$ws.Range("AQ11").Value = 76
$ws.Range("AR11").Value = "female"

# Reflect the reviewer's updated viewport/selection on the sheet
$ws.Activate()
$ws.Range("AL25").Select()
$excel.ActiveWindow.Zoom = 140
